# Apply trade #39 close update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.77   # Current Capital
$wsSummary.Range("B4").Value = -0.23     # Total P&L $
$wsSummary.Range("B5").Value = -0.12     # Total P&L %
$wsSummary.Range("B6").Value = 39        # Total Trades
$wsSummary.Range("B7").Value = 14        # Winning Trades
$wsSummary.Range("B9").Value = 35.9      # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.77
$wsStatus.Range("D4").Value = 39
$wsStatus.Range("E4").Value = -0.23
$wsStatus.Range("F4").Value = -0.23
$wsStatus.Range("G4").Value = 35.9

# --- New trade row (#39), appended to both "All Trades" and "MarketMaking" sheets ---
function Add-TradeRow($ws) {
    $ws.Cells.Item(40, 1).Value = 39
    # Leading apostrophe forces text storage for the date-looking string so
    # Excel does not auto-convert it into a date serial number, matching the
    # plain-text storage used by the rest of the column. The HH:MM:SS time
    # string is stored as text natively, so it needs no such hint.
    $ws.Cells.Item(40, 2).Value = "'2026-02-17"
    $ws.Cells.Item(40, 3).Value = "08:33:08"
    $ws.Cells.Item(40, 4).Value = "MarketMaking"
    $ws.Cells.Item(40, 5).Value = "UP"
    $ws.Cells.Item(40, 6).Value = 0.32
    $ws.Cells.Item(40, 7).Value = 0.43
    $ws.Cells.Item(40, 8).Value = "CLOSED"
    $ws.Cells.Item(40, 9).Value = 34.375
    $ws.Cells.Item(40, 10).Value = 0.11
    $ws.Cells.Item(40, 11).Value = 99.77
    $ws.Cells.Item(40, 12).Value = 0
    $ws.Cells.Item(40, 13).Value = 0
    $ws.Cells.Item(40, 14).Value = 0.6
    $ws.Cells.Item(40, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(40, 16).Value = "early_exit"
    $ws.Cells.Item(40, 17).Value = 0.13
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
